# "Try to use custom message, per user"
# Give row 3 ("Recipients" sheet) its own custom message in column C,
# matching the look (wrap style / row height) of row 2's message cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipients")

# 1) Set the new, recipient-specific message text in C3.
$ws.Range("C3").Value = 'Hi, *Less Important Text*: Lorem Ipsum is simply dummy text of the printing and typesetting industry. Lorem Ipsum has been the industry''s standard dummy text ever since the 1500s, when an unknown printer took a galley of type and scrambled it to make a type specimen book.'

# 2) Row 3 carried a stray leftover row-level format; clear it off the row.
$ws.Rows.Item(3).ClearFormats()

# 3) Re-apply the same per-cell formatting used by the equivalent cells in
#    row 2 (number format on A, wrap-text on C) that ClearFormats removed.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 4) Restore the untouched, style-less spacer cell in column B.
$ws.Range("B1").Copy($ws.Range("B3"))

$excel.CutCopyMode = 0

# 5) Match row 2's (wrapped-text) row height now that C3 has real content.
$ws.Rows.Item(3).RowHeight = 45
